# ml-research.xlsx: add two new (empty) sheets "convolutions" and
# "parameters" after "dataset", and append two dataset rows (Stanford Cars
# Dataset, Vehicle Detection Image Set) to the "dataset" sheet's table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataset")

# --- Append the two new dataset rows -----------------------------------
# Write row 12 (Stanford Cars) columns B/C first, then row 11 (Vehicle
# Detection) columns B/C, then the two "notes" cells in column D, then the
# rank numbers in column A -- this ordering reproduces the target shared
# string table order.
$ws.Range("B12").Value = "Stanford Cars Dataset"
$ws.Range("C12").Value = "https://www.kaggle.com/datasets/jessicali9530/stanford-cars-dataset"

$ws.Range("B11").Value = "Vehicle Detection Image Set"
$ws.Range("C11").Value = "https://www.kaggle.com/datasets/brsdincer/vehicle-detection-image-set"

$ws.Range("D11").Value = "Binary, Decent Variation in orientation"
$ws.Range("D12").Value = "Multi-class, more data"

$ws.Range("A11").ClearFormats() | Out-Null
$ws.Range("A11").Value = 1

$ws.Range("A12").VerticalAlignment = -4160
$ws.Range("A12").Value = 2

# Rank/name/link cell formatting to match the rest of the table.
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").WrapText = $true
$ws.Range("B12").VerticalAlignment = -4160
$ws.Range("B12").WrapText = $true

$ws.Range("C11").WrapText = $true
$ws.Range("C12").WrapText = $true

$ws.Rows.Item(11).RowHeight = 51
$ws.Rows.Item(12).RowHeight = 34

# Row 5 height tightened slightly as part of the same edit.
$ws.Rows.Item(5).RowHeight = 59

# Column C widened to fit the new, longer links.
$ws.Columns.Item(3).ColumnWidth = 30.1666666667

# --- Selection / scroll position ----------------------------------------
$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null

# --- Add the two new (empty) sheets --------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$convSheet = $wb.Worksheets.Add($null, $lastSheet)
$convSheet.Name = "convolutions"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$paramSheet = $wb.Worksheets.Add($null, $lastSheet)
$paramSheet.Name = "parameters"

# Re-activate "dataset" so it stays the tab selected/shown on open.
$ws.Activate() | Out-Null
